# Update performance dashboard 2025-12-19 00:32
# Row 3 (model "gemini-3-pro", Pattern1-Pure Data) changed on both the
# "Summary" sheet and the "Pattern1-Pure Data" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Summary", "Pattern1-Pure Data")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Plain currency-formatted text values - Excel keeps these as text
    # automatically because "\xa5" is not recognized as a numeric currency
    # symbol, so no special handling is required.
    $ws.Range("C3").Value = "¥1,000,000.00"
    $ws.Range("D3").Value = "¥1,001,002.00"
    $ws.Range("E3").Value = "¥+1,002.00"

    # Percent-looking text values - without forcing a text interpretation,
    # Excel would parse these as numeric percentages. Prefix with a quote
    # so they are stored as literal text, matching the source data.
    $ws.Range("F3").Value = "'+0.10%"
    $ws.Range("G3").Value = "'+28.71%"
    $ws.Range("I3").Value = "'0.00%"
    $ws.Range("J3").Value = "'100.0%"
    $ws.Range("K3").Value = "'0.1002%"
}
